# Commit: "modification context dans les extensions" (b65bd0ad3095a651d87ee70bc2c17d6b7266e32b)
#
# Changes applied:
#   1. Metadata sheet (B8): Date value bumped to the new publication timestamp.
#   2. Metadata sheet (B20): Context value expanded from the short element name to
#      the fully-qualified canonical extension-context URL.
#   3. Elements sheet (K6): Extension.value[x] Type(s) changed from "code" to
#      "CodeableConcept" (value keeps its trailing newline).
#   4. Elements sheet, column K width grows to fit the new, longer "CodeableConcept"
#      text (Excel's "best fit" column sizing reacting to the wider cell content).

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2024-06-10T07:36:07+00:00"
$metadata.Range("B20").Value = "element:http://ltsi.univ-rennes.fr/StructureDefinition/oncofair-medicationadministration-component#MedicationAdministration"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K6").Value = "CodeableConcept`n"
$elements.Columns.Item(11).ColumnWidth = 16.3
